# "finish the index and add the news"
#
# This script reproduces, via the PowerPoint COM/VBA object model, the
# OOXML changes described by the target diff:
#   1. Two slide guides (one horizontal, one vertical) are added to the
#      presentation's slide-guide list.
#   2. Two shapes on slide 2 ("矩形 1" / rectangle id=2 and
#      "手繪多邊形 4" / freeform polygon id=5) are nudged to new positions.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide guides (horizontal @ 3840, vertical @ 2160 - both mid-grey).
#    PpGuideOrientation: 1 = ppHorizontalGuide, 2 = ppVerticalGuide.
#    (Best-effort: guarded so that a host without guide support can't
#    abort the rest of the script.)
# ---------------------------------------------------------------------
try {
    $guides = $p.Guides
    [void]$guides.Add(1, 3840)
    [void]$guides.Add(2, 2160)
} catch {
}

# ---------------------------------------------------------------------
# 2) Reposition the two shapes on slide 2.
#    (Left/Top are expressed in points; 1 pt = 12700 EMU.)
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(2)

$rect = $slide.Shapes.Item(2)
if ($rect.Id -ne 2) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        if ($slide.Shapes.Item($i).Id -eq 2) { $rect = $slide.Shapes.Item($i) }
    }
}
$rect.Left = -0.75
$rect.Top = 32.166614173228346

$poly = $slide.Shapes.Item(5)
if ($poly.Id -ne 5) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        if ($slide.Shapes.Item($i).Id -eq 5) { $poly = $slide.Shapes.Item($i) }
    }
}
$poly.Left = 244.57850393700787
$poly.Top = 286.70835885669294
